$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $NewValue
}

# Row 2
Set-TextValue "D2" "290.18"
Set-TextValue "E2" "-4.16%"
# Row 3
Set-TextValue "D3" "30.84"
Set-TextValue "E3" "-4.23%"
# Row 4
Set-TextValue "D4" "4.877"
Set-TextValue "E4" "-2.69%"
# Row 5
Set-TextValue "D5" "0.07178"
Set-TextValue "E5" "-9.25%"
# Row 6
Set-TextValue "D6" "7.688"
Set-TextValue "E6" "-2.35%"
# Row 7
Set-TextValue "D7" "1.741"
Set-TextValue "E7" "-17.42%"
# Row 8
Set-TextValue "D8" "3.739"
Set-TextValue "E8" "-1.90%"
# Row 9
Set-TextValue "E9" "-3.34%"
# Row 10
Set-TextValue "D10" "0.1674"
Set-TextValue "E10" "-4.55%"
# Row 11
Set-TextValue "D11" "0.07323"
Set-TextValue "E11" "-7.68%"
# Row 12
Set-TextValue "D12" "0.08059"
Set-TextValue "E12" "-7.71%"
# Row 13
Set-TextValue "D13" "0.02992"
Set-TextValue "E13" "-4.35%"
# Row 14
Set-TextValue "D14" "0.1000"
Set-TextValue "E14" "-0.31%"
# Row 15
Set-TextValue "D15" "0.001493"
Set-TextValue "E15" "-1.23%"
# Row 16
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04486"
Set-TextValue "E16" "-2.61%"
# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D17" "0.005853"
Set-TextValue "E17" "-1.55%"
# Row 18
$ws.Range("B18").Value = "UpBots"
$ws.Range("C18").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D18" "0.007527"
Set-TextValue "E18" "2,125.44%"
# Row 19
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D19" "3.455"
Set-TextValue "E19" "-0.34%"
# Row 20
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D20" "2.104"
Set-TextValue "E20" "-7.61%"
# Row 21
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D21" "0.3292"
Set-TextValue "E21" "0.14%"
# Row 22
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D22" "0.1302"
Set-TextValue "E22" "0.75%"
# Row 23
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D23" "4.352"
Set-TextValue "E23" "4.46%"
# Row 24
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D24" "0.2003"
Set-TextValue "E24" "11.72%"
# Row 25
Set-TextValue "D25" "0.001213"
Set-TextValue "E25" "-2.01%"
# Row 26
Set-TextValue "D26" "0.004005"
Set-TextValue "E26" "-10.36%"
# Row 27
Set-TextValue "D27" "0.0001253"
Set-TextValue "E27" "0.17%"
# Row 39
Set-TextValue "D39" "0.01639"
Set-TextValue "E39" "-5.07%"
# Row 40
Set-TextValue "D40" "0.04335"
Set-TextValue "E40" "-9.83%"
# Row 41
Set-TextValue "D41" "0.007374"
Set-TextValue "E41" "-0.80%"
# Row 42
Set-TextValue "D42" "0.1310"
Set-TextValue "E42" "-3.94%"
# Row 43
Set-TextValue "D43" "0.002041"
Set-TextValue "E43" "-12.79%"
# Row 44
Set-TextValue "D44" "0.01121"
Set-TextValue "E44" "9.69%"
# Row 45
Set-TextValue "D45" "0.00005729"
Set-TextValue "E45" "-4.45%"
# Row 46
Set-TextValue "E46" "0.17%"
# Row 47
Set-TextValue "D47" "2.187"
Set-TextValue "E47" "165.60%"
# Row 48
Set-TextValue "E48" "-11.26%"
# Row 49
Set-TextValue "D49" "0.00002104"
Set-TextValue "E49" "0.17%"
# Row 50
Set-TextValue "D50" "0.0002004"
